$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "302.94"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "5.36%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "31.86"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "7.77%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "5.223"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.82%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07347"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "8.67%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "7.821"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "6.31%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.735"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "8.52%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "1.485"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "5.23%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9084"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.65%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.01664"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "2,475.52%"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "5.32%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07472"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "7.53%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07996"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "4.01%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.02954"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.09%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09926"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001488"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "-6.40%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.04536"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "1.42%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.006465"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "4.95%"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "3.471"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "0.57%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "2.226"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "-0.20%"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3336"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.28%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1320"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "1.06%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.386"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "7.25%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.1619"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "2.50%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.001218"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "2.16%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.004422"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "6.77%"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.0001299"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "8.48%"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0001740"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04495"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "5.41%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007209"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "5.73%"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "8.04%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002329"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "5.06%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01342"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "2.99%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006076"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "7.67%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-3.57%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.01299"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-13.68%"
